$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dino Breeder 3 - Gaia Fukkatsu: extend comment (G20)
$ws.Range("G20").Value = "Palette 0xE4. Uses a margin of 2 lines instead of 3 usually, you can print the ennemy that you have beaten in a kinf of death march but the game is horribly difficult and confusing."

# Pocket Family 2: fill in previously empty Supported/TX Rate/Compression/Palette cells
# and align their style with the rest of the table (center-aligned, same as column C).
$ws.Range("C76:F76").HorizontalAlignment = -4108
$ws.Range("C76").Value = "Yes"
$ws.Range("D76").Value = "Normal"
$ws.Range("E76").Value = "No"
$ws.Range("F76").Value = "Custom"
$ws.Range("G76").Value = "Palette 0x4C. Original cartridge mandatory, HuC-3 mapper. You can print family pictures with the game title."

# Tales of Phantasia: Nakiri's Dungeon: Supported changes to Partial, comment extended
$ws.Range("C102").Value = "Partial"
$ws.Range("G102").Value = "Palette 0x00, acts as 0xE4 (documented in pandocs). You can print images of your team characters but printing loops infinitely. On real printer it's OK."

# The Little Mermaid 2: Pinball Frenzy: extend comment
$ws.Range("G104").Value = "Palette 0x1B, you can print your scores and some images by playing minigames unlocked on tables."

# VS Lemmings: extend comment
$ws.Range("G108").Value = "Palette 0xE4. Lemmings (USA) has no printer support but it is the same game without VS mode (it is a downgraded version from VS Lemmings)."
